# Request Form - F1.xlsx : "feat: sops Update 4"
#
# Renames the visible worksheet from the old "Software Service Catalog"
# code (S-SW-SC-01) to the "Software Development Lifecycle" code
# (F-SW-SD-01), repoints the Print_Area defined name at the renamed
# sheet, switches the sheet view into Page Layout view scrolled down to
# row 14, and refreshes the footer's revision stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet (S-SW-SC-01 -> F-SW-SD-01) ---------------------
$ws.Name = "F-SW-SD-01"

# --- Repoint the Print_Area defined name at the new sheet name -------
$ws.PageSetup.PrintArea = '$A$1:$F$22'

# --- Switch to Page Layout view, scrolled to row 14 -------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.View = 3            # xlPageLayoutView
$win.ScrollRow = 14      # top-left visible row -> A14
$win.ScrollColumn = 1
$win.TopLeftCell = $ws.Range("A14")

# --- Update the footer revision stamp ---------------------------------
# was: "Rev: 0(0/0/2025)" -> now: "Rev:0(01/10/2025)"
$ws.PageSetup.RightFooter = '&"Arial,Regular"&16Rev:0(01/10/2025)'
